# Commit: "Corrected Lat Lon format or changed spreadsheet tab name from Gliders to Moorings"
#
# Content changes applied:
#   1. Rename worksheet "Gliders" -> "Moorings"
#   2. Update the (now stale) "Glider Serial Number" asset id cell on the renamed
#      sheet from "GP05MOAS-GL340" to "CP05MOAS-GL340"
#   3. Rename the "Glider Serial Number" column header (Asset_Cal_Info sheet) to
#      "Mooring Serial Number"
#   4. Restore the sheet-qualified #REF! defined names that Excel's rename
#      operation otherwise strips the sheet prefix from
#   5. Restore selection/view state for both sheets

$wb = $excel.ActiveWorkbook
$wsMoorings = $wb.Worksheets.Item(1)
$wsAssetCal = $wb.Worksheets.Item(2)

# 1. Rename the "Gliders" tab to "Moorings"
$wsMoorings.Name = "Moorings"

# 2. Rename the "Glider Serial Number" header to "Mooring Serial Number" first so
#    that the new shared string is interned before the asset-id string below
#    (keeps shared-string table ordering consistent with the source edit)
$wsAssetCal.Range("B1").Value = "Mooring Serial Number"

# 3. Correct the serial number value on the Moorings sheet
$wsMoorings.Range("A2").Value = "CP05MOAS-GL340"

# 4. Fix up defined names -- Excel's sheet-rename collapses every "#REF!" defined
#    name (even ones belonging to the sheet that wasn't renamed) down to a bare
#    "#REF!", dropping the sheet qualifier. Restore the proper qualified refs.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase_0").RefersTo = "=Moorings!#REF!"
$wb.Names.Item("_FilterDatabase_0_0_0").RefersTo = "=Moorings!#REF!"
$wb.Names.Item("_FilterDatabase_0_0_0_0_0").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase_0_0_0_0_0_0").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase_0_0_0_0_0_0_0").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase_0_0_0_0_0_0_0_0").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase_0_0_0_0_1").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase_0_0_0_1").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase_0_0_1").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase_0_1").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase_1").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase_1_1").RefersTo = "=Asset_Cal_Info!#REF!"
$wb.Names.Item("_FilterDatabase_2").RefersTo = "=Asset_Cal_Info!#REF!"

# 5. Restore view/selection state: Asset_Cal_Info scrolled back to B1 with B1
#    selected, then return focus to the Moorings tab with D11 selected (the tab
#    that should remain the active one).
$wsAssetCal.Activate() | Out-Null
$wsAssetCal.Range("B1").Select() | Out-Null

$wsMoorings.Activate() | Out-Null
$wsMoorings.Range("D11").Select() | Out-Null
